# Insert a new weekly price record for Naranja / Lane Late (Región Metropolitana)
# at row 15, pushing the existing rows 15-75 down to 16-76.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 15..75 down by one (this also extends the used range to row 76
# and keeps the date-formatted style on column D).
$ws.Rows.Item(15).Insert()

# Populate the newly inserted row 15 with the new record.
$ws.Cells.Item(15, 1).Value = 1
$ws.Cells.Item(15, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(15, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(15, 4).Value = 44525
$ws.Cells.Item(15, 5).Value = 15
$ws.Cells.Item(15, 6).Value = "Fruta"
$ws.Cells.Item(15, 7).Value = 100102
$ws.Cells.Item(15, 8).Value = "Cítricos"
$ws.Cells.Item(15, 9).Value = 100102005
$ws.Cells.Item(15, 10).Value = "Naranja"
$ws.Cells.Item(15, 11).Value = "Lane Late"
$ws.Cells.Item(15, 12).Value = "Segunda"
$ws.Cells.Item(15, 13).Value = 270
$ws.Cells.Item(15, 14).Value = 750
$ws.Cells.Item(15, 15).Value = 800
$ws.Cells.Item(15, 16).Value = 775
$ws.Cells.Item(15, 17).Value = "`$/kilo (en caja de 20 kilos)"
$ws.Cells.Item(15, 18).Value = "Región Metropolitana"
$ws.Cells.Item(15, 19).Value = 775
$ws.Cells.Item(15, 20).Value = 1
